# "Generate Report for Archive"
#
# The report-generation run refreshed the localization status text from
# "Ready for handoff" to "In Translation" wherever it appears (Overview!E2,
# Overview!F2, zh-cn!C2, de-de!C2), which in turn shrank the autosized
# Status/locale columns on the three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Text: "Ready for handoff" -> "In Translation" ------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Column widths: narrow the now-shorter status columns -----------------
# Target stored width is ~13.41 characters; ColumnWidth is quantized to
# whole pixels (MDW=6 for this workbook's default font) on write, so 12.5
# is the input that lands on the closest achievable pixel bucket (13.33).
$narrowWidth = 12.5
$overview.Columns.Item(5).ColumnWidth = $narrowWidth
$overview.Columns.Item(6).ColumnWidth = $narrowWidth
$zhcn.Columns.Item(3).ColumnWidth = $narrowWidth
$dede.Columns.Item(3).ColumnWidth = $narrowWidth
